$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# P1 keeps its text "Split 4" (unchanged content).
# Column R ("Recorded Split 3") is removed entirely, replaced by new headers/data in S:V.
$ws.Range("R1:R6").Clear()
$ws.Range("U1").Value = "Recorded Finish Leg 3"
$ws.Range("V1").Value = "Recorded Finish Leg 4"
$ws.Range("S1").Value = "Recorded Finish Leg 1"
$ws.Range("T1").Value = "Recorded Finish Leg 2"

# --- Move the raw "recorded finish" formulas that used to live in
#     G/J (leg1/leg2) and R/S (leg3/leg4) into the new S/T/U/V block ---
$ws.Range("S2").Formula = '=B16'
$ws.Range("T2").Formula = '=B17'
$ws.Range("U2").Formula = '=B18'
$ws.Range("V2").Formula = '=B19'

$ws.Range("S3").Formula = '=B20'
$ws.Range("T3").Formula = '=B21'
$ws.Range("U3").Formula = '=B22'
$ws.Range("V3").Formula = '=B23'

$ws.Range("S4").Formula = '=B24'
$ws.Range("T4").Formula = '=B25'
$ws.Range("U4").Formula = '=B26'
$ws.Range("V4").Formula = '=B27'

$ws.Range("S5").Formula = '=B28'
$ws.Range("T5").Formula = '=B29'
$ws.Range("U5").Formula = '=B30'
$ws.Range("V5").Formula = '=B31'

$ws.Range("S6").Formula = '=B32'
$ws.Range("T6").Formula = '=B33'
$ws.Range("U6").Formula = '=B34'
$ws.Range("V6").Formula = '=B35'

# --- Rewire the existing leg/split formulas to reference the new S/T/U/V columns ---

# Leg 1 elapsed (F) and split (G) now reference each other via S (recorded finish leg 1)
$ws.Range("F2").Formula = '=S2'
$ws.Range("G2").Formula = '=F2'
$ws.Range("F3").Formula = '=S3-G14'
$ws.Range("G3:G6").Formula = '=F3'
$ws.Range("F5:F6").Formula = '=S5'
$ws.Range("F4").Formula = '=S4'

# Leg 2 elapsed (I) and split (J) reference S/T instead of G/J
$ws.Range("I2:I6").Formula = '=T2-S2'
$ws.Range("J2:J6").Formula = '=G2+I2'

# Leg 3 elapsed (L) references T/U instead of J/R (M - the split - is unchanged)
$ws.Range("L2").Formula = '=U2-MIN(L$10,T2)'
$ws.Range("L3:L6").Formula = '=U3-MIN(L$10,T3)'

# Leg 4 elapsed (O) and split (P) reference U/V instead of R/S
$ws.Range("O2").Formula = '=V2-MIN(O$10,U2)'
$ws.Range("O3").Formula = '=V3-MIN(O$10,U3)'
$ws.Range("O4").Formula = '=V4-MIN(O$10,U4)'
$ws.Range("O5").Formula = '=V5-MIN(O$10,U5)'
$ws.Range("O6").Formula = '=V6-MIN(O$10,U6)'

$ws.Range("P2:P3").Formula = '=M2+O2'
$ws.Range("P4").Formula = '=M4+O4'
$ws.Range("P5:P6").Formula = '=M5+O5'

# --- Update selection to match the edited workbook state ---
$ws.Range("O28").Select()
